# Updated cryptos list refresh (mirrors the GitHub Actions scheduled commit):
# for each affected row in Sheet1, write the new Coin/Link/Price/Volume(1h)
# text values. Every cell in this sheet is plain text (t="inlineStr" in the
# source OOXML), so each assignment below is a literal string. The Price
# column (D) is pre-formatted as Text ("@") before the value is written so
# that numeric-looking text (e.g. "318.70") keeps its exact digits/trailing
# zeros instead of Excel silently auto-converting it to a Number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.839.73'
$ws.Range('E2').Value = '  +1.27%  '

# Row 3 - Ethereum
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.250.30'
$ws.Range('E3').Value = '  +0.48%  '

# Row 4 - TetherUSD
$ws.Range('E4').Value = '  +0.02%  '

# Row 5 - BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.70'
$ws.Range('E5').Value = '  -0.80%  '

# Row 6 - Solana
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.79'
$ws.Range('E6').Value = '  +0.25%  '

# Row 7 - XRP
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.578'
$ws.Range('E7').Value = '  -1.45%  '

# Row 8 - USDC
$ws.Range('E8').Value = '  +0.06%  '

# Row 9 - Cardano
$ws.Range('E9').Value = '  -1.43%  '

# Row 10 - Avalanche
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.06'
$ws.Range('E10').Value = '  -0.11%  '

# Row 11 - Dogecoin
$ws.Range('E11').Value = '  +0.89%  '

# Row 12 - Polkadot
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.72'
$ws.Range('E12').Value = '  -0.05%  '

# Row 13 - TRON
$ws.Range('E13').Value = '  -2.53%  '

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.591.51'
$ws.Range('E14').Value = '  +0.49%  '

# Row 15 - Polygon
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.858'
$ws.Range('E15').Value = '  -1.27%  '

# Row 16 - Chainlink
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.23'
$ws.Range('E16').Value = '  -1.50%  '

# Row 17 - WrappedEther
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.252.22'
$ws.Range('E17').Value = '  +0.06%  '

# Row 18 - WrappedBTC
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.751.66'
$ws.Range('E18').Value = '  +1.17%  '

# Row 19 - InternetComputer(DFINITY)
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.38'
$ws.Range('E19').Value = '  -7.73%  '

# Row 20 - ShibaInu
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0988'
$ws.Range('E20').Value = '  +1.98%  '

# Row 21 - Uniswap
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.57'
$ws.Range('E21').Value = '  +0.47%  '

# Row 22 - Litecoin
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.61'
$ws.Range('E22').Value = '  +0.06%  '

# Row 23 - PancakeSwap
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.14'
$ws.Range('E23').Value = '  -1.60%  '

# Row 24 - BitcoinCash
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '236.17'
$ws.Range('E24').Value = '  -0.72%  '

# Row 25 - ImmutableX
$ws.Range('E25').Value = '  -2.57%  '

# Row 26 - Dai
$ws.Range('E26').Value = '  +0.04%  '

# Row 27 - Cosmos
$ws.Range('E27').Value = '  +0.92%  '

# Row 28 - Toncoin
$ws.Range('E28').Value = '  +0.32%  '

# Row 29 - InjectiveProtocol
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '37.36'
$ws.Range('E29').Value = '  +3.54%  '

# Row 30 - Filecoin
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.25'
$ws.Range('E30').Value = '  -2.20%  '

# Row 31 - Monero
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '159.80'
$ws.Range('E31').Value = '  +4.06%  '

# Row 32 - EthereumClassic
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.18'
$ws.Range('E32').Value = '  -1.38%  '

# Row 33 - Hedera
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0851'
$ws.Range('E33').Value = '  -3.56%  '

# Row 34 - WEMIXToken
$ws.Range('E34').Value = '  -1.22%  '

# Row 35 - Kaspa
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.115'
$ws.Range('E35').Value = '  +9.69%  '

# Row 36 - LidoDAOToken
$ws.Range('E36').Value = '  -6.05%  '

# Row 37 - ARBITRUM
$ws.Range('E37').Value = '  -0.74%  '

# Row 38 - Stellar
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.118'
$ws.Range('E38').Value = '  -3.12%  '

# Row 39 - NEARProtocol
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.77'
$ws.Range('E39').Value = '  +1.77%  '

# Row 40 - RenderToken
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.22'
$ws.Range('E40').Value = '  -5.41%  '

# Row 41 - Celestia
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '15.91'
$ws.Range('E41').Value = '  +17.93%  '

# Row 42 - VeChain
$ws.Range('E42').Value = '  -2.55%  '

# Row 43 - FirstDigitalUSD
$ws.Range('E43').Value = '  +0.15%  '

# Row 44 - Maker
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.798.26'
$ws.Range('E44').Value = '  +1.27%  '

# Rows 46/47 swapped places (ordi <-> BitcoinSV) with updated figures
$ws.Range('B46').Value = 'BitcoinSV'
$ws.Range('C46').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '82.83'
$ws.Range('E46').Value = '  -5.42%  '

$ws.Range('B47').Value = 'ordi'
$ws.Range('C47').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '75.52'
$ws.Range('E47').Value = '  -1.48%  '

# Row 48 - THORChain
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.20'
$ws.Range('E48').Value = '  -2.51%  '

# Row 49 - MultiversX
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '58.85'
$ws.Range('E49').Value = '  -0.86%  '

# Rows 50/51 swapped places (Aave <-> Stacks) with updated figures
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.69'
$ws.Range('E50').Value = '  +5.05%  '

$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '103.72'
$ws.Range('E51').Value = '  -0.15%  '
